$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '20.014.28'
$ws.Range("E2").Value = '  -4.14%  '

$ws.Range("D3").Value = '1.423.29'
$ws.Range("E3").Value = '  -4.34%  '

Set-TextValue $ws.Range("D4") '0.9987'
$ws.Range("E4").Value = '  -0.94%  '

Set-TextValue $ws.Range("D5") '0.9991'
$ws.Range("E5").Value = '  -0.74%  '

Set-TextValue $ws.Range("D6") '277.05'
$ws.Range("E6").Value = '  -1.24%  '

Set-TextValue $ws.Range("D7") '0.3689'
$ws.Range("E7").Value = '  -2.49%  '

Set-TextValue $ws.Range("D8") '0.3111'
$ws.Range("E8").Value = '  +0.61%  '

Set-TextValue $ws.Range("D9") '39.85'
$ws.Range("E9").Value = '  -5.74%  '

$ws.Range("E10").Value = '  +3.02%  '

Set-TextValue $ws.Range("D11") '0.06563'
$ws.Range("E11").Value = '  -3.55%  '

Set-TextValue $ws.Range("D12") '0.9990'
$ws.Range("E12").Value = '  -0.99%  '

Set-TextValue $ws.Range("D13") '5.539'
$ws.Range("E13").Value = '  +1.05%  '

Set-TextValue $ws.Range("D14") '17.80'
$ws.Range("E14").Value = '  +0.62%  '

Set-TextValue $ws.Range("D15") '6.234'
$ws.Range("E15").Value = '  -0.95%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D16") '0.00001027'
$ws.Range("E16").Value = '  -2.38%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.421.18'
$ws.Range("E17").Value = '  -5.14%  '

Set-TextValue $ws.Range("D18") '0.05693'
$ws.Range("E18").Value = '  -12.54%  '

Set-TextValue $ws.Range("D19") '0.9990'
$ws.Range("E19").Value = '  -0.73%  '

Set-TextValue $ws.Range("D20") '71.61'
$ws.Range("E20").Value = '  -11.81%  '

Set-TextValue $ws.Range("D21") '5.644'
$ws.Range("E21").Value = '  -4.29%  '

Set-TextValue $ws.Range("D22") '14.81'
$ws.Range("E22").Value = '  -0.49%  '

Set-TextValue $ws.Range("D23") '11.05'
$ws.Range("E23").Value = '  +2.38%  '

Set-TextValue $ws.Range("D24") '2.239'

$ws.Range("D25").Value = '20.032.98'
$ws.Range("E25").Value = '  -4.10%  '

Set-TextValue $ws.Range("D26") '2.304'
$ws.Range("E26").Value = '  +0.86%  '

Set-TextValue $ws.Range("D27") '133.45'

Set-TextValue $ws.Range("D28") '17.39'
$ws.Range("E28").Value = '  -1.96%  '

$ws.Range("D29").Value = '1.580.08'
$ws.Range("E29").Value = '  -5.11%  '

Set-TextValue $ws.Range("D30") '110.52'
$ws.Range("E30").Value = '  -2.07%  '

Set-TextValue $ws.Range("D31") '3.945'
$ws.Range("E31").Value = '  -17.26%  '

Set-TextValue $ws.Range("D32") '5.304'
$ws.Range("E32").Value = '  -7.26%  '

Set-TextValue $ws.Range("D33") '0.8268'
$ws.Range("E33").Value = '  -9.47%  '

Set-TextValue $ws.Range("D34") '0.07758'
$ws.Range("E34").Value = '  -1.50%  '

Set-TextValue $ws.Range("D35") '1.490'
$ws.Range("E35").Value = '  +1.78%  '

Set-TextValue $ws.Range("D36") '8.329'
$ws.Range("E36").Value = '  -1.47%  '

Set-TextValue $ws.Range("D37") '4.951'
$ws.Range("E37").Value = '  +0.61%  '

Set-TextValue $ws.Range("D38") '0.05861'
$ws.Range("E38").Value = '  +3.38%  '

Set-TextValue $ws.Range("D39") '0.9981'
$ws.Range("E39").Value = '  -0.71%  '

Set-TextValue $ws.Range("D40") '0.02076'
$ws.Range("E40").Value = '  +0.07%  '

Set-TextValue $ws.Range("D41") '10.57'
$ws.Range("E41").Value = '  -3.82%  '

Set-TextValue $ws.Range("D42") '0.1894'
$ws.Range("E42").Value = '  -3.02%  '

Set-TextValue $ws.Range("D43") '1.105'
$ws.Range("E43").Value = '  -2.74%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D44") '0.5346'
$ws.Range("E44").Value = '  -3.13%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D45") '12.42'
$ws.Range("E45").Value = '  -2.81%  '

Set-TextValue $ws.Range("D46") '3.548'
$ws.Range("E46").Value = '  -2.28%  '

Set-TextValue $ws.Range("D47") '0.5227'
$ws.Range("E47").Value = '  -1.88%  '

Set-TextValue $ws.Range("D48") '116.76'
$ws.Range("E48").Value = '  +4.65%  '

Set-TextValue $ws.Range("D49") '1.783'
$ws.Range("E49").Value = '  -1.89%  '

Set-TextValue $ws.Range("D50") '1.039'
$ws.Range("E50").Value = '  -6.73%  '

Set-TextValue $ws.Range("D51") '0.9986'
$ws.Range("E51").Value = '  -0.70%  '
